$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 471
$ws1.Range("F4").Value = 7889
$ws1.Range("F6").Value = 220
$ws1.Range("F8").Value = 30
$ws1.Range("F9").Value = 112
$ws1.Range("F10").Value = 459
$ws1.Range("F15").Value = 71
$ws1.Range("F17").Value = 5788
$ws1.Range("F18").Value = 172
$ws1.Range("F19").Value = 245
$ws1.Range("F20").Value = 1576
$ws1.Range("F22").Value = 361

# Sheet "全部类型" - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 471
$ws4.Range("F4").Value = 7889
$ws4.Range("F6").Value = 220
$ws4.Range("F8").Value = 30
$ws4.Range("F9").Value = 112
$ws4.Range("F10").Value = 459
$ws4.Range("F15").Value = 71
$ws4.Range("F18").Value = 5788
$ws4.Range("F20").Value = 172
$ws4.Range("F21").Value = 245
$ws4.Range("F22").Value = 1576
$ws4.Range("F24").Value = 361
